# Applies the crypto price/volume refresh described in the commit diff.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '27.010.02'
$ws.Range("E2").Value = '  +2.81%  '
$ws.Range("D3").Value = '1.652.49'
$ws.Range("E3").Value = '  +3.64%  '
$ws.Range("E4").Value = '  +0.18%  '
$ws.Range("D5").Formula = '''214.88'
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = '  +1.66%  '
$ws.Range("E6").Value = '  +1.22%  '
$ws.Range("E7").Value = '  +0.11%  '
$ws.Range("E8").Value = '  +1.68%  '
$ws.Range("D9").Formula = '''0.0613'
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = '  +1.59%  '
$ws.Range("D10").Formula = '''19.65'
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = '  +3.42%  '
$ws.Range("E11").Value = '  +1.18%  '
$ws.Range("D12").Value = '1.889.03'
$ws.Range("E12").Value = '  +3.83%  '
$ws.Range("D13").Value = '1.664.09'
$ws.Range("E13").Value = '  +4.30%  '
$ws.Range("E14").Value = '  +2.06%  '
$ws.Range("D15").Formula = '''0.517'
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = '  +2.89%  '
$ws.Range("D16").Formula = '''65.12'
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = '  +2.58%  '
$ws.Range("B17").Value = 'WrappedBTC'
$ws.Range("C17").Value = 'https://coinranking.com/coin/x4WXHge-vvFY+wrappedbtc-wbtc'
$ws.Range("D17").Value = '27.003.31'
$ws.Range("E17").Value = '  +2.81%  '
$ws.Range("B18").Value = 'BitcoinCash'
$ws.Range("C18").Value = 'https://coinranking.com/coin/ZlZpzOJo43mIo+bitcoincash-bch'
$ws.Range("D18").Formula = '''237.95'
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = '  +3.71%  '
$ws.Range("D19").Formula = '''7.81'
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = '  +1.73%  '
$ws.Range("D20").Value = '0.0₃0728'
$ws.Range("E20").Value = '  +1.25%  '
$ws.Range("E21").Value = '  +0.09%  '
$ws.Range("E22").Value = '  +4.36%  '
$ws.Range("E23").Value = '  +4.36%  '
$ws.Range("D24").Formula = '''9.22'
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = '  +3.55%  '
$ws.Range("D25").Formula = '''145.68'
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = '  -0.31%  '
$ws.Range("E26").Value = '  +0.10%  '
$ws.Range("D27").Formula = '''7.10'
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = '  +2.01%  '
$ws.Range("D28").Formula = '''0.113'
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = '  +1.33%  '
$ws.Range("D29").Formula = '''15.80'
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = '  +3.06%  '
$ws.Range("D30").Formula = '''0.0495'
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = '  +0.48%  '
$ws.Range("E31").Value = '  +1.73%  '
$ws.Range("E32").Value = '  +3.15%  '
$ws.Range("D33").Value = '1.511.42'
$ws.Range("E34").Value = '  +4.60%  '
$ws.Range("D35").Formula = '''1.59'
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = '  +8.92%  '
$ws.Range("E36").Value = '  +0.15%  '
$ws.Range("D37").Formula = '''0.574'
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = '  +1.42%  '
$ws.Range("D38").Formula = '''0.885'
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = '  +8.26%  '
$ws.Range("E39").Value = '  +2.68%  '
$ws.Range("D40").Formula = '''5.95'
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = '  +3.56%  '
$ws.Range("E41").Value = '  +0.12%  '
$ws.Range("E42").Value = '  +4.07%  '
$ws.Range("D43").Formula = '''65.79'
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = '  +9.22%  '
$ws.Range("D44").Value = '1.795.46'
$ws.Range("E44").Value = '  +3.60%  '
$ws.Range("D45").Formula = '''0.774'
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = '  +2.71%  '
$ws.Range("D46").Formula = '''0.919'
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = '  -0.88%  '
$ws.Range("D47").Formula = '''89.45'
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = '  +1.33%  '
$ws.Range("E48").Value = '  +0.10%  '
$ws.Range("E49").Value = '  +3.02%  '
$ws.Range("E50").Value = '  +1.09%  '
$ws.Range("D51").Formula = '''0.0974'
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = '  +2.20%  '
